# Insert a new weekly data row for "Femacal de La Calera" / Haba at row 153,
# shifting the existing rows 153-168 down to 154-169.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(153).Insert()

$ws.Cells.Item(153, 1).Value = 3
$ws.Cells.Item(153, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(153, 3).Value = "Coquimbo"
$ws.Cells.Item(153, 4).Value = 44769
$ws.Cells.Item(153, 5).Value = 5
$ws.Cells.Item(153, 6).Value = 100112026
$ws.Cells.Item(153, 7).Value = "Haba"
$ws.Cells.Item(153, 8).Value = "Sin especificar"
$ws.Cells.Item(153, 9).Value = "Primera"
$ws.Cells.Item(153, 10).Value = 105
$ws.Cells.Item(153, 11).Value = 17000
$ws.Cells.Item(153, 12).Value = 18000
$ws.Cells.Item(153, 13).Value = 17524
$ws.Cells.Item(153, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(153, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(153, 16).Value = 701
$ws.Cells.Item(153, 17).Value = 25
$ws.Cells.Item(153, 18).Value = "Hortaliza"

$ws.Range("D153").NumberFormat = $ws.Range("D154").NumberFormat
